$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 140
$ws.Range("H3").Value = 56
$ws.Range("H4").Select()
